$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2403.6
$ws.Range("J58").Value = 3412.1428
$ws.Range("L58").Value = 10236.4284
$ws.Range("N58").Value = -10536.4284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1924.875
$ws.Range("I4").Value = 509.8
$ws.Range("J4").Value = 4283.3335
$ws.Range("K4").Value = 509.8
$ws.Range("L4").Value = 4283.3335
$ws.Range("M4").Value = -393.8
$ws.Range("N4").Value = -4515.3335
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 500
$ws.Range("M29").Value = -192
$ws.Range("H96").Value = 23373
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 23373
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 23373
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -28865
$ws.Range("H132").Value = 1313.9
$ws.Range("I132").Value = 1702
$ws.Range("K132").Value = 5106
$ws.Range("M132").Value = -2576

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1716.3334
$ws.Range("I5").Value = 74.5
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 74.5
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 38.5
$ws.Range("N5").Value = -5226
$ws.Range("H134").Value = 653
$ws.Range("I134").Value = 653
$ws.Range("K134").Value = 1959
$ws.Range("M134").Value = 576

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 826
$ws.Range("I22").Value = 402.42105
$ws.Range("J22").Value = 4850
$ws.Range("K22").Value = 402.42105
$ws.Range("L22").Value = 4850
$ws.Range("M22").Value = -52.42104999999998
$ws.Range("N22").Value = -5550
$ws.Range("H88").Value = 38240.855
$ws.Range("I88").Value = 8000
$ws.Range("J88").Value = 43281
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 43281
$ws.Range("M88").Value = -7594
$ws.Range("N88").Value = -44093
$ws.Range("H91").Value = 38240.855
$ws.Range("I91").Value = 8000
$ws.Range("J91").Value = 43281
$ws.Range("K91").Value = 8000
$ws.Range("L91").Value = 43281
$ws.Range("M91").Value = -6596
$ws.Range("N91").Value = -46089

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H43").Value = 11038.091
$ws.Range("I43").Value = 2700
$ws.Range("J43").Value = 14164.875
$ws.Range("K43").Value = 2700
$ws.Range("L43").Value = 14164.875
$ws.Range("M43").Value = -2549
$ws.Range("N43").Value = -14466.875
$ws.Range("H132").Value = 1579
$ws.Range("I132").Value = 1376
$ws.Range("K132").Value = 4128
$ws.Range("M132").Value = -1598
$ws.Range("H136").Value = 19913.5
$ws.Range("J136").Value = 19913.5
$ws.Range("L136").Value = 59740.5
$ws.Range("N136").Value = -64840.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 40033
$ws.Range("J38").Value = 40033
$ws.Range("L38").Value = 40033
$ws.Range("N38").Value = -40853
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H45").Value = 39999
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H94").Value = 44993
$ws.Range("J94").Value = 44993
$ws.Range("L94").Value = 44993
$ws.Range("N94").Value = -46345
$ws.Range("H124").Value = 45500
$ws.Range("J124").Value = 45500
$ws.Range("L124").Value = 45500
$ws.Range("N124").Value = -55320

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 1000
$ws.Range("M9").Value = -860
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 2000
$ws.Range("J26").Value = 18000
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = -1707
$ws.Range("N26").Value = -18586
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30756
$ws.Range("H80").Value = 31599.2
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 31599.2
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 31599.2
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -33595.2
$ws.Range("H83").Value = 31599.2
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 31599.2
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 94797.6
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -104781.6
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 1756.3846
$ws.Range("I122").Value = 1725.25
$ws.Range("K122").Value = 5175.75
$ws.Range("M122").Value = -2725.75
